# CS449 C tutorial slides - comment updates
$p = $ppt.ActivePresentation

# Slide 1: remove the "Note solutions given in blue are the official solutions."
# textbox (TextBox 4 / shape id 5) that duplicated the note already present
# on slide 2.
$s1 = $p.Slides.Item(1)
for ($i = $s1.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s1.Shapes.Item($i)
    if ($shp.Name -eq "TextBox 4") {
        $shp.Delete()
    }
}

# Slide 4: reword the comment on the return-statement negation.
$s4 = $p.Slides.Item(4)
for ($i = 1; $i -le $s4.Shapes.Count; $i++) {
    $shp = $s4.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -like "Note flip assuming*") {
            $shp.TextFrame.TextRange.Text = "Note negation (assuming two’s complement convention) on return statement."
        }
    }
}
